$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.367.92"
$ws.Range("E2").Value = "  +0.24%  "
Set-TextValue $ws.Range("D3") "3.585.93"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "605.37"
$ws.Range("E5").Value = "  +0.12%  "
Set-TextValue $ws.Range("D6") "148.13"
$ws.Range("E6").Value = "  +3.11%  "
Set-TextValue $ws.Range("D7") "3.584.28"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  -0.27%  "
Set-TextValue $ws.Range("D11") "7.94"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +0.41%  "
Set-TextValue $ws.Range("D13") "4.195.99"
$ws.Range("E13").Value = "  +0.89%  "
Set-TextValue $ws.Range("D14") "0.0000204"
$ws.Range("E14").Value = "  -0.55%  "
Set-TextValue $ws.Range("D15") "29.52"
$ws.Range("E15").Value = "  -1.71%  "
Set-TextValue $ws.Range("D16") "3.585.19"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("E17").Value = "  +1.78%  "
Set-TextValue $ws.Range("D18") "66.410.27"
$ws.Range("E18").Value = "  +0.24%  "
Set-TextValue $ws.Range("D19") "11.07"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("E20").Value = "  +2.40%  "
Set-TextValue $ws.Range("D21") "14.85"
$ws.Range("E21").Value = "  +1.22%  "
Set-TextValue $ws.Range("D22") "422.51"
$ws.Range("E22").Value = "  -1.48%  "
Set-TextValue $ws.Range("D23") "0.610"
$ws.Range("E23").Value = "  +0.31%  "
Set-TextValue $ws.Range("D24") "78.15"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +2.83%  "
Set-TextValue $ws.Range("D29") "2.49"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +0.08%  "
Set-TextValue $ws.Range("D31") "3.584.20"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("E32").Value = "  +3.94%  "
Set-TextValue $ws.Range("D33") "25.00"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  -2.36%  "
Set-TextValue $ws.Range("D36") "7.72"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  -2.45%  "
Set-TextValue $ws.Range("D39") "175.22"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -0.79%  "
Set-TextValue $ws.Range("D43") "46.03"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -3.40%  "
Set-TextValue $ws.Range("D45") "0.999"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  +4.99%  "
Set-TextValue $ws.Range("D47") "23.74"
$ws.Range("E47").Value = "  +3.79%  "
Set-TextValue $ws.Range("D48") "24.24"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("E50").Value = "  -4.66%  "
Set-TextValue $ws.Range("D51") "0.941"
$ws.Range("E51").Value = "  +1.05%  "

# Rows 27 / 28 swapped: RenderToken <-> InternetComputer(DFINITY)
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D27") "9.34"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D28") "8.12"
$ws.Range("E28").Value = "  +3.51%  "
